# Update the Ost (Q) / Nord (R) coordinate cells on row 7 to rounded
# integer values, and clear the Starttid (Z) / Sluttid (AB) time cells
# (Slutdatum in AA stays untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q7").Value = 534450
$ws.Range("R7").Value = 6830575

$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
